{"js": "const pairs = [[\"37+25=\", \"16+63=\"], [\"53+1=\", \"85-29=\"], [\"9+49=\", \"14+75=\"], [\"0+76=\", \"53-18=\"], [\"6+1=\", \"63+6=\"], [\"94-84=\", \"94-42=\"], [\"20-2=\", \"80-8=\"], [\"99-96=\", \"93-10=\"], [\"45+42=\", \"51-24=\"], [\"23+3=\", \"7+70=\"], [\"7-6=\", \"78-45=\"], [\"3+18=\", \"12+66=\"], [\"71-5=\", \"61-48=\"], [\"67+20=\", \"19+59=\"], [\"46+18=\", \"91+8=\"], [\"5+7=\", \"96+3=\"], [\"90-55=\", \"67+28=\"], [\"59+39=\", \"73-6=\"], [\"79-2=\", \"18+37=\"], [\"16+11=\", \"75+7=\"], [\"68-61=\", \"38-11=\"], [\"86-71=\", \"35-1=\"], [\"34+18=\", \"41+52=\"], [\"36+33=\", \"83-19=\"], [\"56-53=\", \"94-66=\"], [\"37-32=\", \"86-63=\"], [\"73-24=\", \"86+4=\"], [\"88-31=\", \"31-0=\"], [\"6+45=\", \"79-32=\"], [\"64-28=\", \"10+51=\"], [\"4+74=\", \"49+34=\"], [\"63+20=\", \"31-17=\"], [\"80-34=\", \"60-29=\"], [\"89-4=\", \"23+70=\"], [\"47+37=\", \"54+16=\"], [\"36+31=\", \"77-60=\"], [\"52+40=\", \"84-41=\"], [\"99-89=\", \"2+91=\"], [\"41+6=\", \"41+6=\"], [\"64-17=\", \"35+43=\"], [\"4-2=\", \"90+4=\"], [\"34+20=\", \"4+71=\"], [\"26+34=\", \"79-9=\"], [\"97-24=\", \"56-36=\"], [\"63-34=\", \"62-0=\"], [\"25+4=\", \"67+14=\"], [\"49-35=\", \"18+5=\"], [\"6+53=\", \"12+58=\"], [\"51+21=\", \"35-19=\"], [\"1+49=\", \"78+11=\"], [\"13+83=\", \"47-4=\"], [\"50-15=\", \"60-6=\"], [\"23-18=\", \"8+34=\"], [\"55+28=\", \"21+70=\"], [\"8+23=\", \"42-13=\"], [\"90-40=\", \"44-5=\"], [\"92-49=\", \"28+13=\"], [\"72-59=\", \"10+46=\"], [\"34+41=\", \"85-59=\"], [\"4+37=\", \"20+18=\"], [\"51+6=\", \"69-6=\"], [\"35-20=\", \"14+48=\"], [\"30+61=\", \"77-61=\"], [\"75+17=\", \"10+44=\"], [\"57-54=\", \"53-2=\"], [\"8+90=\", \"13+85=\"], [\"89-57=\", \"78-24=\"], [\"16+76=\", \"75-32=\"], [\"4+61=\", \"23-6=\"], [\"88-64=\", \"46-37=\"], [\"83-33=\", \"27-2=\"], [\"35-34=\", \"63-9=\"], [\"7+77=\", \"32+22=\"], [\"46-1=\", \"17+48=\"], [\"78-50=\", \"50+29=\"], [\"27+40=\", \"94-93=\"], [\"40-24=\", \"68-14=\"], [\"92-41=\", \"2+8=\"], [\"80-2=\", \"92-31=\"], [\"52-13=\", \"61-37=\"], [\"49+2=\", \"42+7=\"], [\"16+43=\", \"58+26=\"], [\"60+22=\", \"26+5=\"], [\"7+64=\", \"28+27=\"], [\"10+71=\", \"33+54=\"], [\"30-14=\", \"71-63=\"], [\"67+25=\", \"75-26=\"], [\"7+77=\", \"57+37=\"], [\"81-59=\", \"50-17=\"], [\"53-52=\", \"94-21=\"], [\"87-47=\", \"21+10=\"], [\"34-25=\", \"72-13=\"], [\"37+11=\", \"15-2=\"], [\"95-5=\", \"32+44=\"], [\"44+55=\", \"40+58=\"], [\"48+6=\", \"71-38=\"], [\"12+81=\", \"7-7=\"], [\"0+71=\", \"31+3=\"], [\"98-55=\", \"5+27=\"], [\"71-32=\", \"36+42=\"]];\n\n// 1) Update the date heading paragraph (first paragraph in the body).\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst headingPara = paras.items[0];\nconst headingRange = headingPara.getRange();\nheadingRange.load(\"text\");\nawait context.sync();\nif (headingRange.text === \"2023-02-13 Monday\") {\n  headingRange.insertText(\"2023-02-14 Tuesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the arithmetic-problem table, cell by cell, in document order.\n//    A positional mapping is used (rather than a blind global text\n//    replace) because one problem string (\"7+77=\") occurs twice in the\n//    grid with two different replacements.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nlet i = 0;\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    if (i < pairs.length) {\n      const [oldText, newText] = pairs[i];\n      if (values[r][c] === oldText) {\n        values[r][c] = newText;\n      }\n      i++;\n    }\n  }\n}\ntable.values = values;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the body).\n$headingRange = $d.Paragraphs.Item(1).Range\nif ($headingRange.Text.TrimEnd(\"`r\") -eq \"2023-02-13 Monday\") {\n    $headingRange.Text = \"2023-02-14 Tuesday\"\n}\n\n# 2) Update the arithmetic-problem table, cell by cell, in document order.\n#    Expected \"before\" / \"after\" grids are listed explicitly (rather than\n#    doing a blind global find/replace) because one problem string\n#    (\"7+77=\") occurs twice in the table with two different replacements;\n#    only the cell whose current text matches the expected \"before\" value\n#    is touched.\n$oldGrid = @(\n        @(\"37+25=\", \"53+1=\", \"9+49=\", \"0+76=\", \"6+1=\"),\n        @(\"94-84=\", \"20-2=\", \"99-96=\", \"45+42=\", \"23+3=\"),\n        @(\"7-6=\", \"3+18=\", \"71-5=\", \"67+20=\", \"46+18=\"),\n        @(\"5+7=\", \"90-55=\", \"59+39=\", \"79-2=\", \"16+11=\"),\n        @(\"68-61=\", \"86-71=\", \"34+18=\", \"36+33=\", \"56-53=\"),\n        @(\"37-32=\", \"73-24=\", \"88-31=\", \"6+45=\", \"64-28=\"),\n        @(\"4+74=\", \"63+20=\", \"80-34=\", \"89-4=\", \"47+37=\"),\n        @(\"36+31=\", \"52+40=\", \"99-89=\", \"41+6=\", \"64-17=\"),\n        @(\"4-2=\", \"34+20=\", \"26+34=\", \"97-24=\", \"63-34=\"),\n        @(\"25+4=\", \"49-35=\", \"6+53=\", \"51+21=\", \"1+49=\"),\n        @(\"13+83=\", \"50-15=\", \"23-18=\", \"55+28=\", \"8+23=\"),\n        @(\"90-40=\", \"92-49=\", \"72-59=\", \"34+41=\", \"4+37=\"),\n        @(\"51+6=\", \"35-20=\", \"30+61=\", \"75+17=\", \"57-54=\"),\n        @(\"8+90=\", \"89-57=\", \"16+76=\", \"4+61=\", \"88-64=\"),\n        @(\"83-33=\", \"35-34=\", \"7+77=\", \"46-1=\", \"78-50=\"),\n        @(\"27+40=\", \"40-24=\", \"92-41=\", \"80-2=\", \"52-13=\"),\n        @(\"49+2=\", \"16+43=\", \"60+22=\", \"7+64=\", \"10+71=\"),\n        @(\"30-14=\", \"67+25=\", \"7+77=\", \"81-59=\", \"53-52=\"),\n        @(\"87-47=\", \"34-25=\", \"37+11=\", \"95-5=\", \"44+55=\"),\n        @(\"48+6=\", \"12+81=\", \"0+71=\", \"98-55=\", \"71-32=\")\n    )\n\n$newGrid = @(\n        @(\"16+63=\", \"85-29=\", \"14+75=\", \"53-18=\", \"63+6=\"),\n        @(\"94-42=\", \"80-8=\", \"93-10=\", \"51-24=\", \"7+70=\"),\n        @(\"78-45=\", \"12+66=\", \"61-48=\", \"19+59=\", \"91+8=\"),\n        @(\"96+3=\", \"67+28=\", \"73-6=\", \"18+37=\", \"75+7=\"),\n        @(\"38-11=\", \"35-1=\", \"41+52=\", \"83-19=\", \"94-66=\"),\n        @(\"86-63=\", \"86+4=\", \"31-0=\", \"79-32=\", \"10+51=\"),\n        @(\"49+34=\", \"31-17=\", \"60-29=\", \"23+70=\", \"54+16=\"),\n        @(\"77-60=\", \"84-41=\", \"2+91=\", \"41+6=\", \"35+43=\"),\n        @(\"90+4=\", \"4+71=\", \"79-9=\", \"56-36=\", \"62-0=\"),\n        @(\"67+14=\", \"18+5=\", \"12+58=\", \"35-19=\", \"78+11=\"),\n        @(\"47-4=\", \"60-6=\", \"8+34=\", \"21+70=\", \"42-13=\"),\n        @(\"44-5=\", \"28+13=\", \"10+46=\", \"85-59=\", \"20+18=\"),\n        @(\"69-6=\", \"14+48=\", \"77-61=\", \"10+44=\", \"53-2=\"),\n        @(\"13+85=\", \"78-24=\", \"75-32=\", \"23-6=\", \"46-37=\"),\n        @(\"27-2=\", \"63-9=\", \"32+22=\", \"17+48=\", \"50+29=\"),\n        @(\"94-93=\", \"68-14=\", \"2+8=\", \"92-31=\", \"61-37=\"),\n        @(\"42+7=\", \"58+26=\", \"26+5=\", \"28+27=\", \"33+54=\"),\n        @(\"71-63=\", \"75-26=\", \"57+37=\", \"50-17=\", \"94-21=\"),\n        @(\"21+10=\", \"72-13=\", \"15-2=\", \"32+44=\", \"40+58=\"),\n        @(\"71-38=\", \"7-7=\", \"31+3=\", \"5+27=\", \"36+42=\")\n    )\n\n$table = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $oldGrid.Length; $r++) {\n    for ($c = 0; $c -lt $oldGrid[$r].Length; $c++) {\n        $cell = $table.Cell($r + 1, $c + 1)\n        $current = $cell.Range.Text.TrimEnd(\"`r\", [char]7)\n        if ($current -eq $oldGrid[$r][$c]) {\n            $cell.Range.Text = $newGrid[$r][$c]\n        }\n    }\n}\n"}
